# Update "想去人数" (want-to-go count) figures across the workbook sheets
# to reflect the latest generated output (gh-pages build at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 (Exhibitions) ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Cells.Item(2, 6).Value  = 1239
$ws1.Cells.Item(4, 6).Value  = 52
$ws1.Cells.Item(5, 6).Value  = 3497
$ws1.Cells.Item(6, 6).Value  = 1753
$ws1.Cells.Item(7, 6).Value  = 6298
$ws1.Cells.Item(9, 6).Value  = 1890
$ws1.Cells.Item(11, 6).Value = 8
$ws1.Cells.Item(13, 6).Value = 5
$ws1.Cells.Item(14, 6).Value = 4
$ws1.Cells.Item(16, 6).Value = 7484
$ws1.Cells.Item(17, 6).Value = 138
$ws1.Cells.Item(21, 6).Value = 1735
$ws1.Cells.Item(22, 6).Value = 845
$ws1.Cells.Item(28, 6).Value = 1694
$ws1.Cells.Item(29, 6).Value = 791
$ws1.Cells.Item(30, 6).Value = 356
$ws1.Cells.Item(35, 6).Value = 3910

# --- Sheet 2: 演出 (Performances) ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Cells.Item(8, 6).Value = 441

# --- Sheet 3: 本地生活 (Local life) ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Cells.Item(3, 6).Value = 2268
$ws3.Cells.Item(4, 6).Value = 671

# --- Sheet 4: 全部类型 (All types, combined) ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Cells.Item(3, 6).Value  = 2268
$ws4.Cells.Item(4, 6).Value  = 671
$ws4.Cells.Item(5, 6).Value  = 1239
$ws4.Cells.Item(7, 6).Value  = 52
$ws4.Cells.Item(10, 6).Value = 3497
$ws4.Cells.Item(12, 6).Value = 1753
$ws4.Cells.Item(13, 6).Value = 6298
$ws4.Cells.Item(15, 6).Value = 1890
$ws4.Cells.Item(18, 6).Value = 8
$ws4.Cells.Item(20, 6).Value = 5
$ws4.Cells.Item(23, 6).Value = 7484
$ws4.Cells.Item(24, 6).Value = 138
$ws4.Cells.Item(28, 6).Value = 1735
$ws4.Cells.Item(29, 6).Value = 845
$ws4.Cells.Item(34, 6).Value = 1694
$ws4.Cells.Item(35, 6).Value = 791
$ws4.Cells.Item(37, 6).Value = 356
$ws4.Cells.Item(44, 6).Value = 3910
